# Fruta / hortaliza, semanal
# Insert a new week's worth of price data (4 quality grades) at the top of
# the historical block for this market/product, pushing the existing rows
# down by 4 (matching the diff: dimension A1:T738 -> A1:T742).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 635 - everything from row 635 down shifts to 639+
$ws.Rows("635:638").Insert()

# Static column values shared by every row in this market/product block
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100107
$producto    = "Otros"
$categoriaId = 100107011
$categoria   = "Tuna"
$variedad    = "Sin especificar"
$unidad      = "$/caja 18 kilos"
$origen      = "Provincia de Melipilla"
$fecha       = 45212
$kilos       = 18

function Set-TunaRow([int]$row, [string]$calidad, [double]$volumen, [double]$pmin, [double]$pmax, [double]$pmed, [double]$precioKilo) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $pmin
    $ws.Cells.Item($row, 15).Value = $pmax
    $ws.Cells.Item($row, 16).Value = $pmed
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKilo
    $ws.Cells.Item($row, 20).Value = $kilos
}

Set-TunaRow 635 "Especial"               150 30000 30000 30000 1667
Set-TunaRow 636 "Extra (doble especial)" 100 35000 35000 35000 1944
Set-TunaRow 637 "Primera"                 80 25000 25000 25000 1389
Set-TunaRow 638 "Segunda"                 80 20000 20000 20000 1111
